$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first test header row ("Prueba 1") ---
$ws.Range("A2").Value = "pRUEBA 1"
$ws.Range("B2").Value = "2025-09-29_20-53-09"
# C2/D2/E2 stay the same (Epicóndilo lateral / MEDIDAS(Kgf/cm2) / MEDIDAS (kPa))

# --- Rows 3-9: recalibrated measurement values ---
$ws.Range("D3").Value = 4.72
$ws.Range("E3").Value = 463.25

$ws.Range("D4").Value = 4.49
$ws.Range("E4").Value = 440.07

$ws.Range("D5").Value = 4.01
$ws.Range("E5").Value = 393.71

$ws.Range("D6").Value = 4.96
$ws.Range("E6").Value = 486.43

$ws.Range("D7").Value = 4.72
$ws.Range("E7").Value = 463.25

# Row 8 used to hold a second test's header text (A8..E8); now it's
# just two more numeric measurements merged into "Prueba 1".
$ws.Range("A8:C8").ClearContents()
$ws.Range("D8").Value = 4.96
$ws.Range("E8").Value = 486.43

$ws.Range("D9").Value = 4.49
$ws.Range("E9").Value = 440.07

# --- Row 10: becomes the header row for "Prueba 2" ---
$ws.Range("A10").Value = "Prueba 2"
$ws.Range("B10").Value = "2025-09-29_20-54-16"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = "MEDIDAS(Kgf/cm2)"
$ws.Range("E10").Value = "MEDIDAS (kPa)"

# Row 11 used to hold a third test's header text (A11..E11); now it's
# a plain measurement row under "Prueba 2".
$ws.Range("A11:C11").ClearContents()
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 196.16

$ws.Range("D12").Value = 1.88
$ws.Range("E12").Value = 183.89

$ws.Range("D13").Value = 2.25
$ws.Range("E13").Value = 220.71

# --- Remove everything below row 13 (old "Prueba 3" data) ---
$ws.Range("A14:E25").ClearContents()
